$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# The last block on the sheet ("日期：2018.11.15 第十一周四" ... summary) used to
# sit 3 blank rows below the previous block (rows 182-184 empty, block starting
# at row 185). The edit tightens that gap to the normal single blank row by
# removing two of those blank rows, which shifts the whole trailing block
# (old rows 185-193) up to rows 183-191.
$ws.Rows("183:184").Delete()

# Update the block's header date line (old row 185 -> now row 183).
$ws.Range("A183").Value2 = "日期：2018.11.19 第十二周周一"

# 邱志鹏 row (old row 187 -> now row 185): plan text + completion status updated.
$ws.Range("B185").Value2 = "更新完善与后台的对接"
$ws.Range("C185").Value2 = "完成"

# 李达波 row (old row 190 -> now row 188): plan text + completion status updated.
$ws.Range("B188").Value2 = "完成APP端和后台对接工作，编写网页端界面。"
$ws.Range("C188").Value2 = "完成"

# 黄俊贤 row (old row 189 -> now row 187): plan text updated (was on leave).
$ws.Range("B187").Value2 = "完善头像选择功能"

# 冯德志 row (old row 191 -> now row 189): plan text shortened.
$ws.Range("B189").Value2 = "整合两份地图代码"
